$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.481.94"
$ws.Range("E2").Value = "  +1.00%  "

# Row 3
$ws.Range("D3").Value = "1.919.63"
$ws.Range("E3").Value = "  +1.60%  "

# Row 4
$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = "  +0.83%  "

# Row 5
$ws.Range("D5").Value = "'325.09"
$ws.Range("E5").Value = "  +0.92%  "

# Row 6
$ws.Range("E6").Value = "  +0.64%  "

# Row 7
$ws.Range("D7").Value = "'0.4833"
$ws.Range("E7").Value = "  +2.77%  "

# Row 8
$ws.Range("D8").Value = "'0.4086"
$ws.Range("E8").Value = "  +1.59%  "

# Row 9
$ws.Range("D9").Value = "'0.08173"
$ws.Range("E9").Value = "  +2.10%  "

# Row 10
$ws.Range("D10").Value = "'1.024"
$ws.Range("E10").Value = "  +3.03%  "

# Row 11
$ws.Range("D11").Value = "'23.56"
$ws.Range("E11").Value = "  +4.29%  "

# Row 12
$ws.Range("D12").Value = "1.898.25"
$ws.Range("E12").Value = "  +2.17%  "

# Row 13
$ws.Range("D13").Value = "'6.042"
$ws.Range("E13").Value = "  +2.71%  "

# Row 14
$ws.Range("D14").Value = "'7.241"
$ws.Range("E14").Value = "  +2.99%  "

# Row 15
$ws.Range("D15").Value = "'91.30"
$ws.Range("E15").Value = "  +2.74%  "

# Row 16
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").Value = "'0.06790"
$ws.Range("E16").Value = "  +2.62%  "

# Row 17
$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").Value = "'1.007"
$ws.Range("E17").Value = "  +0.59%  "

# Row 18
$ws.Range("E18").Value = "  +1.53%  "

# Row 19
$ws.Range("E19").Value = "  +1.89%  "

# Row 20
$ws.Range("E20").Value = "  +0.63%  "

# Row 21
$ws.Range("D21").Value = "29.505.50"
$ws.Range("E21").Value = "  +1.06%  "

# Row 22
$ws.Range("D22").Value = "'5.638"
$ws.Range("E22").Value = "  +2.71%  "

# Row 23
$ws.Range("D23").Value = "'11.77"
$ws.Range("E23").Value = "  +0.83%  "

# Row 24
$ws.Range("D24").Value = "'2.183"
$ws.Range("E24").Value = "  +0.27%  "

# Row 25
$ws.Range("D25").Value = "2.146.66"
$ws.Range("E25").Value = "  -0.87%  "

# Row 26
$ws.Range("D26").Value = "'6.771"
$ws.Range("E26").Value = "  +12.32%  "

# Row 27
$ws.Range("D27").Value = "'157.27"
$ws.Range("E27").Value = "  +1.51%  "

# Row 28
$ws.Range("D28").Value = "'20.05"
$ws.Range("E28").Value = "  +2.08%  "

# Row 29
$ws.Range("D29").Value = "'2.119"
$ws.Range("E29").Value = "  +1.96%  "

# Row 30
$ws.Range("D30").Value = "'120.60"
$ws.Range("E30").Value = "  +2.16%  "

# Row 31
$ws.Range("D31").Value = "'1.028"
$ws.Range("E31").Value = "  -0.11%  "

# Row 32
$ws.Range("D32").Value = "'0.09576"
$ws.Range("E32").Value = "  +1.72%  "

# Row 33
$ws.Range("D33").Value = "'5.530"
$ws.Range("E33").Value = "  +3.31%  "

# Row 34
$ws.Range("D34").Value = "'3.569"
$ws.Range("E34").Value = "  +0.79%  "

# Row 35
$ws.Range("D35").Value = "'1.388"
$ws.Range("E35").Value = "  +0.46%  "

# Row 36
$ws.Range("D36").Value = "'0.02282"
$ws.Range("E36").Value = "  +2.42%  "

# Row 37
$ws.Range("D37").Value = "'0.06146"
$ws.Range("E37").Value = "  +1.32%  "

# Row 38
$ws.Range("D38").Value = "'1.179"
$ws.Range("E38").Value = "  +0.64%  "

# Row 39
$ws.Range("D39").Value = "'0.5986"
$ws.Range("E39").Value = "  +3.03%  "

# Row 40
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").Value = "'10.85"
$ws.Range("E40").Value = "  +8.40%  "

# Row 41
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'8.042"
$ws.Range("E41").Value = "  +0.29%  "

# Row 42
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.1861"
$ws.Range("E42").Value = "  +1.86%  "

# Row 43
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'2.426"
$ws.Range("E43").Value = "  -1.60%  "

# Row 44
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "'1.284"
$ws.Range("E44").Value = "  +0.81%  "

# Row 45
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "'0.07606"
$ws.Range("E45").Value = "  -1.46%  "

# Row 46
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'12.41"
$ws.Range("E46").Value = "  +2.61%  "

# Row 47
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.5585"
$ws.Range("E47").Value = "  +2.15%  "

# Row 48
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'1.964"
$ws.Range("E48").Value = "  +3.42%  "

# Row 49
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'117.70"
$ws.Range("E49").Value = "  +3.72%  "

# Row 50
$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").Value = "'2.431"
$ws.Range("E50").Value = "  +4.27%  "

# Row 51
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'72.74"
$ws.Range("E51").Value = "  +2.60%  "

